$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1; 3 = 1; 4 = 0; 5 = 0; 7 = 0; 8 = 1; 9 = 2; 10 = 1; 11 = 1; 12 = 0;
    13 = 0; 14 = 3; 15 = 0; 16 = 0; 17 = 1; 18 = 0; 19 = 0; 20 = 2; 21 = 0; 22 = 1;
    23 = 1; 24 = 0; 25 = 2; 26 = 0; 27 = 0; 28 = 0; 29 = 1; 31 = 1; 32 = 0; 33 = 0;
    34 = 0; 35 = 1; 36 = 2; 38 = 1; 39 = 0; 40 = 2; 41 = 0; 42 = 1; 43 = 0; 44 = 0;
    45 = 0; 46 = 2; 47 = 1; 48 = 1; 49 = 0; 50 = 1; 51 = 0; 52 = 0; 53 = 0; 54 = 2;
    55 = 1; 56 = 1; 57 = 3; 58 = 1; 59 = 1; 60 = 0; 61 = 2; 62 = 2; 63 = 1; 65 = 0;
    66 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
